$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 582.5769
$ws.Range("I15").Value = 582.5769
$ws.Range("K15").Value = 1747.7307
$ws.Range("M15").Value = -1578.7307

$ws.Range("H20").Value = 4704.1665
$ws.Range("I20").Value = 645
$ws.Range("J20").Value = 25000
$ws.Range("K20").Value = 645
$ws.Range("L20").Value = 25000
$ws.Range("M20").Value = -415
$ws.Range("N20").Value = -25460

$ws.Range("H35").Value = 4704.1665
$ws.Range("I35").Value = 645
$ws.Range("J35").Value = 25000
$ws.Range("K35").Value = 645
$ws.Range("L35").Value = 25000
$ws.Range("M35").Value = -266
$ws.Range("N35").Value = -25758

$ws.Range("H42").Value = 14.333333
$ws.Range("I42").Value = 14.333333
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 42.999999
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 187.000001
$ws.Range("N42").ClearContents()

$ws.Range("H96").Value = 2080.1667
$ws.Range("I96").Value = 370.5
$ws.Range("K96").Value = 1111.5
$ws.Range("M96").Value = 261.5

$ws.Range("H132").Value = 5227.7
$ws.Range("I132").Value = 2046.8334
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 6140.5002
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -3610.5002
$ws.Range("N132").Value = -35057

$ws.Range("H137").Value = 1265.12
$ws.Range("I137").Value = 1142.2273
$ws.Range("J137").Value = 2166.3333
$ws.Range("K137").Value = 3426.6819
$ws.Range("L137").Value = 6498.999899999999
$ws.Range("M137").Value = -876.6819
$ws.Range("N137").Value = -11598.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7035.5557
$ws.Range("I61").Value = 7035.5557
$ws.Range("K61").Value = 7035.5557
$ws.Range("M61").Value = -6823.5557

$ws.Range("H102").Value = 2591.9167
$ws.Range("I102").Value = 1138
$ws.Range("J102").Value = 5499.75
$ws.Range("K102").Value = 1138
$ws.Range("L102").Value = 5499.75
$ws.Range("M102").Value = 484
$ws.Range("N102").Value = -8743.75

$ws.Range("H136").Value = 7035.5557
$ws.Range("I136").Value = 7035.5557
$ws.Range("K136").Value = 21106.6671
$ws.Range("M136").Value = -18556.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 5000
$ws.Range("J19").Value = 5000
$ws.Range("L19").Value = 5000
$ws.Range("N19").Value = -5346

$ws.Range("H82").Value = 11375
$ws.Range("I82").Value = 11375
$ws.Range("K82").Value = 11375
$ws.Range("M82").Value = -10992

$ws.Range("H85").Value = 11375
$ws.Range("I85").Value = 11375
$ws.Range("K85").Value = 11375
$ws.Range("M85").Value = -10049

$ws.Range("H94").Value = 1340.7
$ws.Range("J94").Value = 1598.2858
$ws.Range("L94").Value = 1598.2858
$ws.Range("N94").Value = -2500.2858

$ws.Range("H107").Value = 596.36365
$ws.Range("I107").Value = 596.36365
$ws.Range("K107").Value = 596.36365
$ws.Range("M107").Value = 1323.63635

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 959.8
$ws.Range("I16").Value = 949.75
$ws.Range("K16").Value = 949.75
$ws.Range("M16").Value = -662.75

$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 2000
$ws.Range("N22").Value = -2700
$ws.Range("M22").ClearContents()

$ws.Range("H31").Value = 4487.4375
$ws.Range("I31").Value = 2982.3333
$ws.Range("K31").Value = 2982.3333
$ws.Range("M31").Value = -2687.3333

$ws.Range("H34").Value = 4487.4375
$ws.Range("I34").Value = 2982.3333
$ws.Range("K34").Value = 2982.3333
$ws.Range("M34").Value = -2780.3333

$ws.Range("H107").Value = 459.55
$ws.Range("I107").Value = 340.25
$ws.Range("J107").Value = 936.75
$ws.Range("K107").Value = 340.25
$ws.Range("L107").Value = 936.75
$ws.Range("M107").Value = 1579.75
$ws.Range("N107").Value = -4776.75

$ws.Range("H113").Value = 959.8
$ws.Range("I113").Value = 949.75
$ws.Range("K113").Value = 949.75
$ws.Range("M113").Value = 1220.25

$ws.Range("H132").Value = 6117.609
$ws.Range("I132").Value = 5285.25
$ws.Range("K132").Value = 15855.75
$ws.Range("M132").Value = -13325.75

$ws.Range("H134").Value = 2475.077
$ws.Range("J134").Value = 2097.5
$ws.Range("L134").Value = 6292.5
$ws.Range("N134").Value = -11362.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1000
$ws.Range("J92").Value = 1000
$ws.Range("L92").Value = 3000
$ws.Range("N92").Value = -5496

$ws.Range("H97").Value = 700
$ws.Range("J97").Value = 400
$ws.Range("L97").Value = 1200
$ws.Range("N97").Value = -2192

$ws.Range("H98").Value = 571.875
$ws.Range("J98").Value = 649.6667
$ws.Range("L98").Value = 1949.0001
$ws.Range("N98").Value = -4945.0001

$ws.Range("H103").Value = 1485.2858
$ws.Range("I103").Value = 342.25
$ws.Range("J103").Value = 3009.3333
$ws.Range("K103").Value = 1026.75
$ws.Range("L103").Value = 9027.999899999999
$ws.Range("M103").Value = -147.75
$ws.Range("N103").Value = -10785.9999

$ws.Range("H109").Value = 558.5
$ws.Range("I109").Value = 409.4
$ws.Range("J109").Value = 807
$ws.Range("K109").Value = 1228.2
$ws.Range("L109").Value = 2421
$ws.Range("M109").Value = -188.1999999999998
$ws.Range("N109").Value = -4501

$ws.Range("H134").Value = 3249
$ws.Range("I134").Value = 998.6667
$ws.Range("K134").Value = 2996.0001
$ws.Range("M134").Value = 2073.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

$ws.Range("H46").Value = 11547.625
$ws.Range("I46").Value = 4166.6665
$ws.Range("J46").Value = 15976.2
$ws.Range("K46").Value = 4166.6665
$ws.Range("L46").Value = 15976.2
$ws.Range("M46").Value = -4010.6665
$ws.Range("N46").Value = -16288.2

$ws.Range("H132").Value = 3287.625
$ws.Range("I132").Value = 2050.5
$ws.Range("K132").Value = 6151.5
$ws.Range("M132").Value = -3621.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4563.273
$ws.Range("J82").Value = 4919.6
$ws.Range("L82").Value = 4919.6
$ws.Range("N82").Value = -5641.6

$ws.Range("H85").Value = 4563.273
$ws.Range("J85").Value = 4919.6
$ws.Range("L85").Value = 4919.6
$ws.Range("N85").Value = -7415.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H81").Value = 2087.5
$ws.Range("J81").Value = 2283
$ws.Range("L81").Value = 4566
$ws.Range("N81").Value = -6688

$ws.Range("H84").Value = 2087.5
$ws.Range("J84").Value = 2283
$ws.Range("L84").Value = 22830
$ws.Range("N84").Value = -33438

$ws.Range("H107").Value = 478.5
$ws.Range("I107").Value = 479.2
$ws.Range("J107").Value = 477.33334
$ws.Range("K107").Value = 1437.6
$ws.Range("L107").Value = 1432.00002
$ws.Range("M107").Value = 482.4000000000001
$ws.Range("N107").Value = -5272.000019999999

$ws.Range("H136").Value = 2673.2666
$ws.Range("I136").Value = 2435.6428
$ws.Range("K136").Value = 7306.928400000001
$ws.Range("M136").Value = -4756.928400000001
